# Generate Report for Handback
# Row 7 in both the zh-cn and de-de sheets corresponds to the file
# 326272d6-7e03-43f0-b301-7a6561ed92a8.md. A new handback xliff was
# detected for that file in both locales, so the "Latest Target File",
# "Latest Handback File" and "Latest Handback DateTime" columns (I/J/K)
# get populated, and an "Error Detail" (P) is recorded because the
# handback was produced against a stale commit of the source file.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/26f4a644e41d06dc77d359da049753cd48a35d31/e2e/326272d6-7e03-43f0-b301-7a6561ed92a8.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f2b9c2adf5d5297537a1b774f98cdc99a995c929/e2e/326272d6-7e03-43f0-b301-7a6561ed92a8.md."
$targetUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f2b9c2adf5d5297537a1b774f98cdc99a995c929/e2e/326272d6-7e03-43f0-b301-7a6561ed92a8.md"

$sheetNames = @("zh-cn", "de-de")
$handbackDateTimes = @("2016-08-25 20:55:20", "2016-08-25 20:55:27")
$handbackFiles = @(
    "326272d6-7e03-43f0-b301-7a6561ed92a8.0ca41bb57beb9d6526b765136eb8c3a473c6a4d1.zh-cn.xlf",
    "326272d6-7e03-43f0-b301-7a6561ed92a8.0ca41bb57beb9d6526b765136eb8c3a473c6a4d1.de-de.xlf"
)

for ($i = 0; $i -lt $sheetNames.Count; $i++) {
    $ws = $wb.Worksheets.Item($sheetNames[$i])

    $ws.Range("I7").Value = "326272d6-7e03-43f0-b301-7a6561ed92a8.md"
    $ws.Hyperlinks.Add($ws.Range("I7"), $targetUrl, "", "", "326272d6-7e03-43f0-b301-7a6561ed92a8.md") | Out-Null

    $ws.Range("J7").Value = $handbackFiles[$i]
    $ws.Range("K7").Value = $handbackDateTimes[$i]
    $ws.Range("P7").Value = $errorDetail
}
